# Edit script implementing the commit "File updates from RMI through 12/4"
# for InputData/fuels/BS/BAU Subsidies.xlsx
#
# Summary of semantic changes (everything else in the underlying XML diff is
# a mechanical side effect of the shared-string table being compacted after
# one string is deleted and a new one appended -- Excel/COM handles that
# automatically as soon as the cell text itself is correct):
#   1. Calculations!B4 and Calculations!B10 label text changes from
#      "Model output, due to endogenous learning" to "See elec/CCaMC".
#   2. Calculations!C4:AH4 and Calculations!C10:AH10 raw data updated to new
#      figures (all downstream formulas / other sheets recompute from this).
#   3. A handful of cosmetic sheet selection changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2: Calculations sheet label + data updates
# ---------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calculations")

$wsCalc.Range("B4").Value = "See elec/CCaMC"
$wsCalc.Range("B10").Value = "See elec/CCaMC"

$row4Values = @(
    1462790.2009077901,
    1284828.0428722501,
    1237544.2340474708,
    1190260.4252226916,
    1142976.6163979124,
    1095692.8075731331,
    1048408.9987483537,
    1001125.1899235747,
    953841.38109879522,
    906557.57227401587,
    859273.76344923663,
    811989.95462445728,
    764706.14579967817,
    757913.21727074252,
    751120.28874180699,
    744327.36021287122,
    737534.43168393557,
    730741.50315499993,
    723948.57462606428,
    717155.64609712863,
    710362.71756819298,
    703569.78903925733,
    696776.86051032168,
    689983.93198138592,
    683191.00345245027,
    676398.07492351462,
    669605.14639457897,
    662812.21786564332,
    656019.28933670768,
    649226.36080777203,
    642433.43227883638,
    635640.50374990073
)

$row10Values = @(
    6831836.4795198459,
    6500515.8491650894,
    6169195.2188103329,
    5843882.7444435377,
    5630238.7966989167,
    5422794.6290817559,
    5233450.6532794116,
    5058216.8964543967,
    4898609.074266511,
    4753976.9239266422,
    4620219.8524136795,
    4501029.6773458235,
    4393594.5403691512,
    4297491.1596871642,
    4212533.3580054678,
    4135494.7609418505,
    4069499.8884631144,
    4010225.5494396384,
    3958592.5068339193,
    3914721.5043390612,
    3876446.5635061474,
    3843367.6747004823,
    3814239.6224356443,
    3789350.4384214408,
    3768365.5466244677,
    3749025.879948609,
    3730662.2863270584,
    3714585.4148475262,
    3697197.8942802823,
    3680144.9497464434,
    3662115.9321582974,
    3640824.119373824
)

# Columns C..AH = column index 3..34
$col = 3
foreach ($val in $row4Values) {
    $wsCalc.Cells.Item(4, $col).Value = $val
    $col = $col + 1
}

$col = 3
foreach ($val in $row10Values) {
    $wsCalc.Cells.Item(10, $col).Value = $val
    $col = $col + 1
}

# ---------------------------------------------------------------------
# 3: Cosmetic selection changes
# ---------------------------------------------------------------------
$wsSubsidies = $wb.Worksheets.Item("Subsidies Paid")
$wsSubsidies.Range("A8:XFD8").Select()

$wsCalc.Range("B21").Select()

$wsBSfTFpEUP = $wb.Worksheets.Item("BS-BSfTFpEUP")
$wsBSfTFpEUP.Range("AG5").Select()

$wsBSpUEO = $wb.Worksheets.Item("BS-BSpUEO")
$wsBSpUEO.Range("B6").Select()

# Restore the originally active sheet/selection (About!B41) so the workbook
# re-opens on the same tab it started on.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B41").Select()
